# Re-run of the nikon_cellpose_bags_spots analysis with updated
# condition labels ("P_None_No" / "None_0" replacing "lig1ON_P" /
# "allON_HP" for the re-bucketed ROIs) and additional images worth of
# per-ROI spot rows, ahead of cleaning up the verbose console output.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("All_Data")

# Rows 9-24 on All_Data: overwrite existing rows 9-17 with new values,
# and append new rows 18-24 (dimension grows from J17 to J24).
$ws1.Range("A9").Value = 1
$ws1.Range("B9").Value = 120
$ws1.Range("C9").Value = 535894
$ws1.Range("D9").Value = 27920.17667460005
$ws1.Range("E9").Value = 558403.5334920011
$ws1.Range("F9").Value = 0.004297967072291778
$ws1.Range("G9").Value = 0.0002148983536145889
$ws1.Range("H9").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I9").Value = "P_None_No"
$ws1.Range("J9").Value = 2

$ws1.Range("A10").Value = 2
$ws1.Range("B10").Value = 101
$ws1.Range("C10").Value = 482984
$ws1.Range("D10").Value = 25163.5558729992
$ws1.Range("E10").Value = 503271.1174599839
$ws1.Range("F10").Value = 0.004013741162407584
$ws1.Range("G10").Value = 0.0002006870581203792
$ws1.Range("H10").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I10").Value = "P_None_No"
$ws1.Range("J10").Value = 2

$ws1.Range("A11").Value = 2
$ws1.Range("B11").Value = 114
$ws1.Range("C11").Value = 565643
$ws1.Range("D11").Value = 29470.10508561544
$ws1.Range("E11").Value = 589402.1017123087
$ws1.Range("F11").Value = 0.003868326891567285
$ws1.Range("G11").Value = 0.0001934163445783643
$ws1.Range("H11").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I11").Value = "P_None_No"
$ws1.Range("J11").Value = 3

$ws1.Range("A12").Value = 1
$ws1.Range("B12").Value = 91
$ws1.Range("C12").Value = 267342
$ws1.Range("D12").Value = 13928.56772522351
$ws1.Range("E12").Value = 278571.3545044702
$ws1.Range("F12").Value = 0.006533335070425537
$ws1.Range("G12").Value = 0.0003266667535212769
$ws1.Range("H12").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I12").Value = "P_None_No"
$ws1.Range("J12").Value = 4

$ws1.Range("A13").Value = 2
$ws1.Range("B13").Value = 108
$ws1.Range("C13").Value = 281652
$ws1.Range("D13").Value = 14674.12137615733
$ws1.Range("E13").Value = 293482.4275231465
$ws1.Range("F13").Value = 0.007359895507984525
$ws1.Range("G13").Value = 0.0003679947753992263
$ws1.Range("H13").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I13").Value = "P_None_No"
$ws1.Range("J13").Value = 4

$ws1.Range("A14").Value = 3
$ws1.Range("B14").Value = 122
$ws1.Range("C14").Value = 269392
$ws1.Range("D14").Value = 14035.37310498691
$ws1.Range("E14").Value = 280707.4620997383
$ws1.Range("F14").Value = 0.008692323252643146
$ws1.Range("G14").Value = 0.0004346161626321573
$ws1.Range("H14").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I14").Value = "P_None_No"
$ws1.Range("J14").Value = 4

$ws1.Range("A15").Value = 3
$ws1.Range("B15").Value = 114
$ws1.Range("C15").Value = 416163
$ws1.Range("D15").Value = 21682.16939437945
$ws1.Range("E15").Value = 433643.3878875891
$ws1.Range("F15").Value = 0.005257776467217879
$ws1.Range("G15").Value = 0.0002628888233608939
$ws1.Range("H15").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I15").Value = "P_None_No"
$ws1.Range("J15").Value = 6

$ws1.Range("A16").Value = 3
$ws1.Range("B16").Value = 115
$ws1.Range("C16").Value = 481085
$ws1.Range("D16").Value = 25064.61762120861
$ws1.Range("E16").Value = 501292.3524241722
$ws1.Range("F16").Value = 0.004588141009687373
$ws1.Range("G16").Value = 0.0002294070504843687
$ws1.Range("H16").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I16").Value = "P_None_No"
$ws1.Range("J16").Value = 7

$ws1.Range("A17").Value = 4
$ws1.Range("B17").Value = 123
$ws1.Range("C17").Value = 454033
$ws1.Range("D17").Value = 23655.20340981367
$ws1.Range("E17").Value = 473104.0681962734
$ws1.Range("F17").Value = 0.005199701641498963
$ws1.Range("G17").Value = 0.0002599850820749481
$ws1.Range("H17").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I17").Value = "P_None_No"
$ws1.Range("J17").Value = 7

$ws1.Range("A18").Value = 5
$ws1.Range("B18").Value = 130
$ws1.Range("C18").Value = 439263
$ws1.Range("D18").Value = 22885.68367366465
$ws1.Range("E18").Value = 457713.673473293
$ws1.Range("F18").Value = 0.005680407098766095
$ws1.Range("G18").Value = 0.0002840203549383048
$ws1.Range("H18").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I18").Value = "P_None_No"
$ws1.Range("J18").Value = 7

$ws1.Range("A19").Value = 1
$ws1.Range("B19").Value = 122
$ws1.Range("C19").Value = 750346
$ws1.Range("D19").Value = 39093.16560192771
$ws1.Range("E19").Value = 781863.3120385543
$ws1.Range("F19").Value = 0.003120750088194037
$ws1.Range("G19").Value = 0.0001560375044097018
$ws1.Range("H19").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I19").Value = "P_None_No"
$ws1.Range("J19").Value = 8

$ws1.Range("A20").Value = 2
$ws1.Range("B20").Value = 131
$ws1.Range("C20").Value = 717257
$ws1.Range("D20").Value = 37369.22257217586
$ws1.Range("E20").Value = 747384.4514435171
$ws1.Range("F20").Value = 0.003505558611688624
$ws1.Range("G20").Value = 0.0001752779305844312
$ws1.Range("H20").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I20").Value = "P_None_No"
$ws1.Range("J20").Value = 8

$ws1.Range("A21").Value = 1
$ws1.Range("B21").Value = 121
$ws1.Range("C21").Value = 375962
$ws1.Range("D21").Value = 19587.68984712646
$ws1.Range("E21").Value = 391753.7969425292
$ws1.Range("F21").Value = 0.006177349189432407
$ws1.Range("G21").Value = 0.0003088674594716203
$ws1.Range("H21").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I21").Value = "P_None_No"
$ws1.Range("J21").Value = 9

$ws1.Range("A22").Value = 2
$ws1.Range("B22").Value = 107
$ws1.Range("C22").Value = 348514
$ws1.Range("D22").Value = 18157.64396237234
$ws1.Range("E22").Value = 363152.8792474468
$ws1.Range("F22").Value = 0.005892835007764971
$ws1.Range("G22").Value = 0.0002946417503882486
$ws1.Range("H22").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I22").Value = "P_None_No"
$ws1.Range("J22").Value = 9

$ws1.Range("A23").Value = 1
$ws1.Range("B23").Value = 7
$ws1.Range("C23").Value = 286682
$ws1.Range("D23").Value = 14936.18530796705
$ws1.Range("E23").Value = 298723.7061593409
$ws1.Range("F23").Value = 0.0004686604950104736
$ws1.Range("G23").Value = 0.00002343302475052368
$ws1.Range("H23").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I23").Value = "None_0"
$ws1.Range("J23").Value = 1

$ws1.Range("A24").Value = 1
$ws1.Range("B24").Value = 6
$ws1.Range("C24").Value = 944338
$ws1.Range("D24").Value = 49200.18473903135
$ws1.Range("E24").Value = 984003.694780627
$ws1.Range("F24").Value = 0.0001219507616043583
$ws1.Range("G24").Value = 0.000006097538080217915
$ws1.Range("H24").Value = "nikon_cellpose_bags_spots"
$ws1.Range("I24").Value = "None_0"
$ws1.Range("J24").Value = 5

$ws2 = $wb.Worksheets.Item("Summary_by_Condition")

# Summary_by_Condition: row 2 becomes None_0, row 3 becomes P_None_No,
# and a new row 4 holds the data that used to be row 2 (allON_HP).
$ws2.Range("A2").Value = "None_0"
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 6.5
$ws2.Range("D2").Value = 0.707
$ws2.Range("E2").Value = 13
$ws2.Range("F2").Value = 32068.185
$ws2.Range("G2").Value = 24228.306
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 0
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0

$ws2.Range("A3").Value = "P_None_No"
$ws2.Range("B3").Value = 14
$ws2.Range("C3").Value = 115.643
$ws2.Range("D3").Value = 11.043
$ws2.Range("E3").Value = 1619
$ws2.Range("F3").Value = 23763.378
$ws2.Range("G3").Value = 7852.325
$ws2.Range("H3").Value = 0.005
$ws2.Range("I3").Value = 0.002
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0

$ws2.Range("A4").Value = "allON_HP"
$ws2.Range("B4").Value = 7
$ws2.Range("C4").Value = 650.429
$ws2.Range("D4").Value = 115.299
$ws2.Range("E4").Value = 4553
$ws2.Range("F4").Value = 22036.741
$ws2.Range("G4").Value = 5712.431
$ws2.Range("H4").Value = 0.03
$ws2.Range("I4").Value = 0.003
$ws2.Range("J4").Value = 0.002
$ws2.Range("K4").Value = 0
